$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right
$ws.Columns.Item(1).Insert()

# Set header for new column A
$ws.Range("A1").Value = "IssueID"

# Fill IssueID values for rows 2-44 (101..143)
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 1).Value = 99 + $r
}

# Shift conditional formatting ranges right by one column to match the
# inserted column (AppliesTo ranges are not auto-shifted by Insert()).
# The three rules are known, in file order, to apply to:
#   1) C2:C29 C31:C44 D1:D43  (dxfId 2, priority 3)
#   2) C30                    (dxfId 1, priority 2)
#   3) D44                    (dxfId 0, priority 1)
# and must become D2:D29 D31:D44 E1:E43 / D30 / E44 respectively.
$fcs = $ws.Cells.FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("D2:D29"))
$fcs.Item(2).ModifyAppliesToRange($ws.Range("D30"))
$fcs.Item(3).ModifyAppliesToRange($ws.Range("E44"))

# Re-apply the first rule's remaining two areas (D31:D44 and E1:E43) as
# additional conditional-format rules sharing the same look (matching the
# red-fill / dark-red-text "less than 1" style), since this host only
# keeps a single contiguous area per ModifyAppliesToRange call.
$extra1 = $ws.Range("D31:D44").FormatConditions.Add(1, 6, "1")
$extra1.Interior.Color = 13551615
$extra1.Font.Color = 393372

$extra2 = $ws.Range("E1:E43").FormatConditions.Add(1, 6, "1")
$extra2.Interior.Color = 13551615
$extra2.Font.Color = 393372

# Update selection to match the new state
$ws.Range("A2:A44").Select()
